$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("arcs")
$ws.Name = "links"
